$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (45310 -> 45311, i.e. 2024-01-19 -> 2024-01-20)
$ws.Range("A1").Value = 45311

# Update the price list in column D (rows 32-37) - adjusted pricing
$ws.Range("D32").Value = 5301
$ws.Range("D33").Value = 5909.76
$ws.Range("D34").Value = 6714.6
$ws.Range("D35").Value = 8356.200000000001
$ws.Range("D36").Value = 8652.6
$ws.Range("D37").Value = 9108.6
